$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowVals($ws, $r) {
    return @(
        $ws.Range("A$r").Value2,
        $ws.Range("B$r").Value2,
        $ws.Range("C$r").Value2,
        $ws.Range("D$r").Value2,
        $ws.Range("E$r").Value2
    )
}

function Set-RowVals($ws, $r, $vals) {
    # Column A is a non-numeric rating label - assign directly.
    $ws.Range("A$r").Value2 = $vals[0]
    # Columns B-E hold numeric-looking text (e.g. "0.8447", "208") that must stay
    # text, not become real numbers - prefix with an apostrophe to force text entry.
    $ws.Range("B$r").Value2 = "'" + $vals[1]
    $ws.Range("C$r").Value2 = "'" + $vals[2]
    $ws.Range("D$r").Value2 = "'" + $vals[3]
    $ws.Range("E$r").Value2 = "'" + $vals[4]
}

# Fix classification report sort order: AAA, AA, A, BBB, BB, B, CCC, CC, C, D
# Swap row 2 (A) and row 4 (AAA)
$row2 = Get-RowVals $ws 2
$row4 = Get-RowVals $ws 4
Set-RowVals $ws 2 $row4
Set-RowVals $ws 4 $row2

# Swap row 5 (B) and row 7 (BBB)
$row5 = Get-RowVals $ws 5
$row7 = Get-RowVals $ws 7
Set-RowVals $ws 5 $row7
Set-RowVals $ws 7 $row5

# Swap row 8 (C) and row 10 (CCC)
$row8 = Get-RowVals $ws 8
$row10 = Get-RowVals $ws 10
Set-RowVals $ws 8 $row10
Set-RowVals $ws 10 $row8
